$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Main")

# --- Sector label: "W" -> "Biopharma" (C5), and drop the now-redundant
#     hyperlink display override (the linked friendly name already reads
#     "Biopharma") ---
$ws.Range("C5").Value = "Biopharma"

$c5Link = $null
foreach ($h in $ws.Hyperlinks) {
    if ($h.Range.Address() -eq '$C$5') {
        $c5Link = $h
    }
}
if ($c5Link -ne $null) {
    $c5Link.TextToDisplay = $null
}

# --- New date cell J3, styled like the other "small font" dates (D5/D11/...) ---
$ws.Range("D5").Copy()
$ws.Range("J3").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("J3").Value = 44989

# --- New date cell G18, default font + short-date number format ---
$ws.Range("G18").Value = 44989
$ws.Range("G18").NumberFormat = "mm-dd-yy"

# --- Refreshed dates ---
$ws.Range("D5").Value = 44987
$ws.Range("D11").Value = 44987
$ws.Range("C20").Value = 44989

# --- Selection moved from I6 to D6 ---
$ws.Range("D6").Select()
